$wb = $excel.ActiveWorkbook

# The "汽車" (car) sheet is the 3rd sheet in the workbook.
$ws = $wb.Worksheets.Item(3)

# ------------------------------------------------------------------
# Row 1: header row
# ------------------------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# New header cells H1:N1 - give them the same (bold/centered/bordered) look
# as the rest of the header row by copying G1's style first.
$ws.Range("G1").Copy($ws.Range("H1:N1"))
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ------------------------------------------------------------------
# Row 2: data row
# ------------------------------------------------------------------
# B2, C2, D2, E2, F2, G2 already hold the correct values (toyotarav4rod,
# 2362, 林淑芬, 100年03月01曰, 買賣, 989000) - nothing to change there.

# New data cells H2:N2 - match the plain data style used by the rest of
# the row by copying G2's style first.
$ws.Range("G2").Copy($ws.Range("H2:N2"))
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# J2 ("date" column) must hold the literal text "2012-03-01" rather than
# being auto-converted into a date serial number. Build it as a formula
# that evaluates to that text in a scratch cell, then paste only the
# resulting value into J2 so it lands as a plain string cell (matching the
# rest of the sheet) without pulling in the formula or changing its style.
$stage = $ws.Range("ZZ1")
$stage.Formula = '="2012-03-01"'
$stage.Copy()
$ws.Range("J2").PasteSpecial(-4163)
$stage.Clear()

$ws.Range("K2").Value = "林淑芬"
$ws.Range("L2").Value = 1337
$ws.Range("M2").Value = "tmp3f851"
$ws.Range("N2").Value = 30

"done"
